# example_month_2.xlsx — shift the logged dates from the old placeholder
# year (Jan/Feb 2001) to Feb/Mar 2025, and blank out the days that no
# longer have entries once the month is re-based (Feb 2025 only has 28
# days, so what used to be rows 30-34 loses its Discharge/Sticker data,
# and the trailing template rows 32-39 go fully blank in column A too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: re-date rows 2..31 (old serials 36892..36921 -> new
#     serials 45689..45718, i.e. 2025-02-01 .. 2025-03-02) ---
$newStart = 45689
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 1).Value = $newStart + ($row - 2)
}

# --- Rows 30 & 31 keep their (new) date but lose the Discharge (B) and
#     Sticker (D) entries that used to live there ---
$ws.Range("B30:D31").Clear()

# --- Rows 32..34 lose their date entirely along with B/D, while keeping
#     the existing cell formatting on column A (same as rows 35-39) ---
$ws.Range("B32:D34").Clear()
$ws.Range("A32:A34").ClearContents()

# --- Update the view: scrolled down to row 21, active cell now C33 ---
$ws.Range("A21").Select()
$ws.Range("C33").Select()

